# Add "Italy" and "Spain" test-data sheets (copies of the "Norway" sheet
# layout/format) to the end of the workbook, fill in their market-specific
# values, and leave "Italy" as the active sheet/tab (mirrors the author's
# "Added test data for ItalyFC Market" commit).

$wb = $excel.ActiveWorkbook
$norway = $wb.Worksheets.Item("Norway")

# --- Italy sheet: clone Norway's layout, rename, fill in values ---
$norway.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3443/T1968"

# --- Spain sheet: clone Norway's layout, rename, fill in values ---
$norway.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3442/T2127"

# Set each new sheet's remembered selection (order matters: last
# activated/selected sheet becomes the workbook's active tab).
$spain.Activate()
$spain.Range("B7").Select()

$italy.Activate()
$italy.Range("A9").Select()
